$d = $word.ActiveDocument

# The "Requisitos" section ends with a paragraph mentioning
# "(Requisito fraco)". It used to be followed by:
#   - a blank paragraph
#   - "Ver no Jupiter Salvar em pdf Salvar em docx"
#   - "© 2020 . Contact: ... Powered by Jekyll and Github pages. ..."
# before the remaining blank paragraph / page-break paragraph at the
# end of the document. Those three paragraphs (the site-chrome
# boilerplate) are removed, leaving the trailing blank paragraph and
# the page-break paragraph untouched.

$startPara = $null
$endPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*Requisito fraco*") {
        $startPara = $i
    }
    if ($t -like "*Jekyll and Github pages*") {
        $endPara = $i
    }
}

if ($startPara -ne $null -and $endPara -ne $null -and $endPara -gt $startPara) {
    $rangeStart = $d.Paragraphs.Item($startPara + 1).Range.Start
    $rangeEnd = $d.Paragraphs.Item($endPara).Range.End

    $r = $d.Range($rangeStart, $rangeEnd)
    $r.Delete()
}
